# Insert a new data row at row 147 (pushing the existing rows 147-211 down
# to 148-212, dimension grows from A1:R211 to A1:R212), then populate the
# newly inserted row with its own record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(147).Insert()

$ws.Cells.Item(147, 1).Value  = 8
$ws.Cells.Item(147, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(147, 3).Value  = "Coquimbo"
$ws.Cells.Item(147, 4).Value  = 44489
$ws.Cells.Item(147, 5).Value  = 4
$ws.Cells.Item(147, 6).Value  = 100114013
$ws.Cells.Item(147, 7).Value  = "Zanahoria"
$ws.Cells.Item(147, 8).Value  = "Sin especificar"
$ws.Cells.Item(147, 9).Value  = "Primera"
$ws.Cells.Item(147, 10).Value = 800
$ws.Cells.Item(147, 11).Value = 6500
$ws.Cells.Item(147, 12).Value = 7000
$ws.Cells.Item(147, 13).Value = 6750
$ws.Cells.Item(147, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(147, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(147, 16).Value = 338
$ws.Cells.Item(147, 17).Value = 20
$ws.Cells.Item(147, 18).Value = "Hortaliza"
